# Changed font on Figures/Tables captions to 9pt.
#
# Each caption textbox ("Fig N - ..." / "Table N - ...") goes from 10pt to
# 9pt. Because the textboxes auto-fit to their text (a:spAutoFit), shrinking
# the font also shrinks/repositions the box slightly, so we update each
# shape's position/size to match the new auto-fit geometry as well.

$p = $ppt.ActivePresentation

$EmuPerPt = 12700.0

function Set-CaptionBox($SlideIndex, $ShapeName, $LeftPt, $TopPt, $WidthPt, $HeightPt) {
    $slide = $p.Slides.Item($SlideIndex)
    $shape = $slide.Shapes.Item($ShapeName)

    # Shrink the caption text from 10pt to 9pt.
    $shape.TextFrame.TextRange.Font.Size = 9

    # Update the auto-fit box geometry to match the new (smaller) text size.
    $shape.Left = $LeftPt
    $shape.Top = $TopPt
    $shape.Width = $WidthPt
    $shape.Height = $HeightPt
}

# Common width/height shared by all the caption boxes (EMU -> pt).
$w = 3619500 / $EmuPerPt
$h = 230832 / $EmuPerPt

# Slide 5 - Fig 1 caption (position unchanged).
$s5Left = 542925 / $EmuPerPt
$s5Top = 6267450 / $EmuPerPt
Set-CaptionBox 5 "TextBox 4" $s5Left $s5Top $w $h

# Slide 6 - Fig 2 caption.
$s6Left = 552450 / $EmuPerPt
$s6Top = 5953125 / $EmuPerPt
Set-CaptionBox 6 "TextBox 4" $s6Left $s6Top $w $h

# Slide 7 - Fig 3 caption.
$s7Left = 542925 / $EmuPerPt
$s7Top = 5972175 / $EmuPerPt
Set-CaptionBox 7 "TextBox 4" $s7Left $s7Top $w $h

# Slide 8 - Table 1 caption.
$s8t1Left = 523875 / $EmuPerPt
$s8t1Top = 4686300 / $EmuPerPt
Set-CaptionBox 8 "TextBox 4" $s8t1Left $s8t1Top $w $h

# Slide 8 - Table 2 caption.
$s8t2Left = 523875 / $EmuPerPt
$s8t2Top = 6334125 / $EmuPerPt
Set-CaptionBox 8 "TextBox 10" $s8t2Left $s8t2Top $w $h

# Slide 9 - Table 3 caption.
$s9Left = 561975 / $EmuPerPt
$s9Top = 5524500 / $EmuPerPt
Set-CaptionBox 9 "TextBox 4" $s9Left $s9Top $w $h
